# Apply the "consolidated report" attendance fixes to the sheet.
# Column H (Absent) values are corrected / filled in for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("H15").Value = 0
